# Adds a "01/01/2023" block of 7 rows (one per equipment type, incl. Total)
# for each of the three regions (Brasil, Nordeste, Sergipe) immediately
# after that region's existing "01/01/2022" block, matching the upstream
# dataset update described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$equip = @("Mamógrafo", "Raio X", "Tomógrafo Computadorizado", "Ressonância Magnética", "Ultrassom", "Equipo Odontológico Completo", "Total")

# Values keyed by region, in the same order as $equip
$brasil2023   = @(2.323793599616808, 37.57924715273676, 3.042754810799485, 1.541952179610437, 23.99203614757576, 77.50078208415611, 145.9805659744954)
$nordeste2023 = @(1.99918668826086, 25.30585624126762, 2.110825120165461, 1.006463401478406, 20.31647709230044, 60.3946741460523, 111.1334826895251)
$sergipe2023  = @(1.809423900469903, 18.72543338858388, 1.430707270138993, 0.7574332606618199, 18.220477881476, 48.72820643591042, 89.67168213724102)

# Process regions from bottom to top so earlier insert points are unaffected
# by later (already-performed) insertions.
# Original (pre-edit) last row of each region's 2012-2022 data block:
#   Brasil:   row 78
#   Nordeste: row 155
#   Sergipe:  row 232

$blocks = @(
    @{ StartRow = 233; Region = "Sergipe";  Values = $sergipe2023 },
    @{ StartRow = 156; Region = "Nordeste"; Values = $nordeste2023 },
    @{ StartRow = 79;  Region = "Brasil";   Values = $brasil2023 }
)

foreach ($block in $blocks) {
    $startRow = $block.StartRow
    $endRow = $startRow + 6

    # Insert 7 fresh rows, pushing everything from $startRow down.
    $ws.Range("A" + $startRow + ":A" + $endRow).EntireRow.Insert()

    for ($i = 0; $i -lt 7; $i++) {
        $r = $startRow + $i
        $ws.Cells.Item($r, 1).Value = $block.Region
        $ws.Cells.Item($r, 2).Value = $equip[$i]
        # Force the date-like text to stay a string rather than being
        # auto-converted into a date serial number.
        $ws.Cells.Item($r, 3).Value = "'01/01/2023"
        $ws.Cells.Item($r, 4).Value = $block.Values[$i]
    }
}
